$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Header row: swap labels - A1 becomes "maturities", B1 becomes "quotes"
$ws.Range("A1").Value = "maturities"
$ws.Range("B1").Value = "quotes"

# Maturity labels in column A become text like "12M" instead of numeric 12
$ws.Range("A2").Value = "12M"
$ws.Range("A3").Value = "24M"
$ws.Range("A4").Value = "36M"
$ws.Range("A5").Value = "72M"
$ws.Range("A6").Value = "120M"
$ws.Range("A7").Value = "240M"

# Update selection to A8
$ws.Range("A8").Select()
